# Applies the cryptos list refresh described in the commit:
# "Updated cryptos list on Mon Dec  2 13:33:54 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.731.76"
$ws.Range("E2").Value = "  -1.62%  "

$ws.Range("D3").Value = "3.615.32"
$ws.Range("E3").Value = "  -2.36%  "

$ws.Range("B4").Value = "XRP"
$ws.Range("C4").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2.41"
$ws.Range("E4").Value = "  +25.94%  "

$ws.Range("B5").Value = "TetherUSD"
$ws.Range("C5").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.00"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "225.22"
$ws.Range("E6").Value = "  -5.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "638.33"
$ws.Range("E7").Value = "  -3.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.419"
$ws.Range("E8").Value = "  -1.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.10"
$ws.Range("E9").Value = "  +3.16%  "

$ws.Range("D11").Value = "3.612.67"
$ws.Range("E11").Value = "  -2.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.40"
$ws.Range("E12").Value = "  +9.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.207"
$ws.Range("E13").Value = "  -0.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000290"
$ws.Range("E14").Value = "  -6.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.53"
$ws.Range("E15").Value = "  -3.71%  "

$ws.Range("D16").Value = "4.289.98"
$ws.Range("E16").Value = "  -2.34%  "

$ws.Range("D17").Value = "95.327.39"
$ws.Range("E17").Value = "  -1.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.81"
$ws.Range("E18").Value = "  -4.08%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.89"
$ws.Range("E19").Value = "  +11.54%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.06"
$ws.Range("E20").Value = "  +8.04%  "

$ws.Range("D21").Value = "3.615.13"
$ws.Range("E21").Value = "  -2.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.523"
$ws.Range("E22").Value = "  +3.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "510.45"
$ws.Range("E23").Value = "  -2.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.24"
$ws.Range("E24").Value = "  -5.94%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "123.20"
$ws.Range("E25").Value = "  +21.02%  "

$ws.Range("B26").Value = "Hedera"
$ws.Range("C26").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.241"
$ws.Range("E26").Value = "  +22.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000202"
$ws.Range("E27").Value = "  -4.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.78"
$ws.Range("E28").Value = "  -1.98%  "

$ws.Range("D29").Value = "3.808.87"
$ws.Range("E29").Value = "  -2.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.68"
$ws.Range("E30").Value = "  -6.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.23"
$ws.Range("E31").Value = "  +4.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.98"
$ws.Range("E32").Value = "  -2.77%  "

$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.181"
$ws.Range("E34").Value = "  -5.27%  "

$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.614"
$ws.Range("E36").Value = "  +2.87%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "32.30"
$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").Value = "  -6.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "591.42"
$ws.Range("E40").Value = "  -9.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.39"
$ws.Range("E41").Value = "  -5.64%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.04"
$ws.Range("E42").Value = "  +3.13%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.39"
$ws.Range("E43").Value = "  +4.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.160"
$ws.Range("E44").Value = "  -2.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.485"
$ws.Range("E45").Value = "  -1.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0487"
$ws.Range("E46").Value = "  +5.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.94"
$ws.Range("E47").Value = "  -6.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.932"
$ws.Range("E48").Value = "  -4.13%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.29"
$ws.Range("E49").Value = "  -0.94%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.70"
$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.50"
$ws.Range("E51").Value = "  -0.67%  "
